$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Ruwayd"
$ws.Range("B5").Value = "ESC"
$ws.Range("C5").Value = "rmushtaq"
$ws.Range("D5").Value = "UC"
$ws.Range("E5").Value = "Compuet"
